# Insert a new row above row 2 on the "stimuli" sheet, shifting the
# existing randomisation table down by one row, and set the new row's
# PID to 1999. Also make the "stimuli" sheet the active/selected tab
# (instead of "stimulus frequency count").

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("stimuli")

# Insert a new row above row 2 (pushes rows 2:42 down to 3:43, and
# expands the bound table from A1:F42 to A1:F43).
$ws.Rows.Item(2).Insert()

# Set the PID for the newly inserted row.
$ws.Range("A2").Value = 1999

# Make the "stimuli" sheet the active sheet/tab.
$ws.Activate()
